$d = $word.ActiveDocument

# The skills bullet list ends with "..., MinIO object store, Java." and the
# edit appends ", SQL" right after "Java" (before the final period), as its
# own run so that the surrounding run structure/formatting is preserved:
#   ...<w:r><w:t>Java</w:t></w:r><w:r><w:t>, SQL</w:t></w:r><w:r><w:t>.</w:t></w:r>

$range = $d.Content
$found = $range.Find.Execute("Java", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)

if ($found) {
    $javaStart = $range.Start
    $javaEnd = $range.End

    # Insert the new text right after "Java" (and before the following
    # "." run), keeping the same character formatting (sz=18 / szCs=18,
    # i.e. 9pt) as the "Java" run it is attached to.
    $range.Collapse(0)
    $range.InsertAfter(", SQL")

    # Force Word to keep the newly inserted text as its own run instead of
    # silently re-merging it into the neighbouring identically formatted
    # runs: toggle Bold off/on/off so a real formatting boundary is
    # recorded, then restore the original (unbolded) look.
    $range.Bold = 1
    $range.Bold = 0

    # Do the same for the original "Java" run so it doesn't stay merged
    # with the run that precedes it (", ") either.
    $javaRange = $d.Range($javaStart, $javaEnd)
    $javaRange.Bold = 1
    $javaRange.Bold = 0
}
